$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Rows.Item(98).Insert(-4121)
$ws.Range("A99:Q99").Copy()
$ws.Range("A98:Q98").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A98:B98").Merge()
$ws.Range("C98:G98").Merge()
$ws.Range("H98:K98").Merge()
$ws.Range("L98:M98").Merge()
$ws.Range("N98:O98").Merge()
Write-Host "done merges"
